$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 13 (blank label in column A, with the "5840535 - Messias Borges
# Silva" text in B/C) is removed entirely; everything below shifts up by one
# row (row heights come along for free with the native row delete).
$ws.Rows.Item(13).Delete()

# After the shift, a handful of cells carry content that differs from what a
# plain shift would have produced - patch those in directly.
$ws.Range("B10").Value = "5840535 - Messias Borges Silva"
$ws.Range("C10").Value = "5840535 - Messias Borges Silva"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# "01/01/2018" looks like a date, and a plain .Value assignment would get
# smart-parsed into a date serial (changing both the stored value and the
# cell's number format/style). Stage the literal text in a scratch cell
# forced to Text format, copy it, and paste-special *values only* into the
# destination cells so B15/C15 keep their original style/format untouched.
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "01/01/2018"
$ws.Range("Z1").Copy()
$ws.Range("B15").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("C15").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("Z1").Clear()

$ws.Range("B18").Value = "5840535 - Messias Borges Silva"
$ws.Range("C18").Value = "5840535 - Messias Borges Silva"

$ws.Range("B19").Value = "Aulas Expositivas; trabalhos e seminários."
$ws.Range("C19").Value = "Aulas Expositivas; trabalhos e seminários."

$ws.Range("B20").Value = "MF = (0,30*P1 + 0,30*P2 + 0,40*TRAB), onde P1 e P2 são provas e TRAB é a nota média de trabalhos e seminários."
$ws.Range("C20").Value = "MF = (0,30*P1 + 0,30*P2 + 0,40*TRAB), onde P1 e P2 são provas e TRAB é a nota média de trabalhos e seminários."

$ws.Range("B21").Value = "NF = (MF + PR)/2, onde PR é uma prova de recuperação."
$ws.Range("C21").Value = "NF = (MF + PR)/2, onde PR é uma prova de recuperação."
